$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9159.799999999999
$ws.Range("I32").Value = 1933.3334
$ws.Range("J32").Value = 19999.5
$ws.Range("K32").Value = 1933.3334
$ws.Range("L32").Value = 19999.5
$ws.Range("M32").Value = -1607.3334
$ws.Range("N32").Value = -20651.5

$ws.Range("H111").Value = 739.8889
$ws.Range("I111").Value = 692.2
$ws.Range("J111").Value = 799.5
$ws.Range("K111").Value = 2076.6
$ws.Range("L111").Value = 2398.5
$ws.Range("M111").Value = 990.3999999999996
$ws.Range("N111").Value = -8532.5

$ws.Range("H112").Value = 3060.6274
$ws.Range("I112").Value = 1332
$ws.Range("J112").Value = 3168.6667
$ws.Range("K112").Value = 3996
$ws.Range("L112").Value = 9506.000100000001
$ws.Range("M112").Value = -2888
$ws.Range("N112").Value = -11722.0001

$ws.Range("H115").Value = 7089.6665
$ws.Range("I115").Value = 1274.5
$ws.Range("J115").Value = 9997.25
$ws.Range("K115").Value = 3823.5
$ws.Range("L115").Value = 29991.75
$ws.Range("M115").Value = -2256.5
$ws.Range("N115").Value = -33125.75

$ws.Range("H118").Value = 1721.8462
$ws.Range("I118").Value = 1657
$ws.Range("J118").Value = 2500
$ws.Range("K118").Value = 4971
$ws.Range("L118").Value = 7500
$ws.Range("M118").Value = -3314
$ws.Range("N118").Value = -10814

$ws.Range("H125").Value = 15499.25
$ws.Range("I125").Value = 17332.666
$ws.Range("J125").Value = 9999
$ws.Range("K125").Value = 155993.994
$ws.Range("L125").Value = 89991
$ws.Range("M125").Value = -153533.994
$ws.Range("N125").Value = -94911

$ws.Range("H132").Value = 2837.463
$ws.Range("I132").Value = 2464.8774
$ws.Range("J132").Value = 6488.8
$ws.Range("K132").Value = 7394.6322
$ws.Range("L132").Value = 19466.4
$ws.Range("M132").Value = -4864.6322
$ws.Range("N132").Value = -24526.4

$ws.Range("H137").Value = 1513.5957
$ws.Range("I137").Value = 1252.7567
$ws.Range("J137").Value = 2478.7
$ws.Range("K137").Value = 3758.2701
$ws.Range("L137").Value = 7436.099999999999
$ws.Range("M137").Value = -1208.2701
$ws.Range("N137").Value = -12536.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 344801.5
$ws.Range("I8").Value = 671670
$ws.Range("J8").Value = 17933
$ws.Range("K8").Value = 671670
$ws.Range("L8").Value = 17933
$ws.Range("M8").Value = -671526
$ws.Range("N8").Value = -18221

$ws.Range("H11").Value = 12509038
$ws.Range("I11").Value = 16667434
$ws.Range("J11").Value = 10014000
$ws.Range("K11").Value = 16667434
$ws.Range("L11").Value = 10014000
$ws.Range("M11").Value = -16667290
$ws.Range("N11").Value = -10014288

$ws.Range("H13").Value = 8475.75
$ws.Range("I13").Value = 1951.5
$ws.Range("J13").Value = 15000
$ws.Range("K13").Value = 1951.5
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = -1807.5
$ws.Range("N13").Value = -15288

$ws.Range("H32").Value = 7928.242
$ws.Range("I32").Value = 1794.6666
$ws.Range("J32").Value = 18662
$ws.Range("K32").Value = 1794.6666
$ws.Range("L32").Value = 18662
$ws.Range("M32").Value = -1507.6666
$ws.Range("N32").Value = -19236

$ws.Range("H74").Value = 3047.8333
$ws.Range("I74").Value = 2585.6875
$ws.Range("J74").Value = 6745
$ws.Range("K74").Value = 2585.6875
$ws.Range("L74").Value = 6745
$ws.Range("M74").Value = -1711.6875
$ws.Range("N74").Value = -8493

$ws.Range("H77").Value = 3047.8333
$ws.Range("I77").Value = 2585.6875
$ws.Range("J77").Value = 6745
$ws.Range("K77").Value = 12928.4375
$ws.Range("L77").Value = 33725
$ws.Range("M77").Value = -8560.4375
$ws.Range("N77").Value = -42461

$ws.Range("H132").Value = 2245.1191
$ws.Range("I132").Value = 1877.3549
$ws.Range("J132").Value = 3281.5454
$ws.Range("K132").Value = 5632.0647
$ws.Range("L132").Value = 9844.636200000001
$ws.Range("M132").Value = -3102.0647
$ws.Range("N132").Value = -14904.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5360.4
$ws.Range("I5").Value = 2950.5
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 2950.5
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -2837.5
$ws.Range("N5").Value = -15226

$ws.Range("H7").Value = 5054125
$ws.Range("I7").Value = 3849.4
$ws.Range("J7").Value = 10104401
$ws.Range("K7").Value = 3849.4
$ws.Range("L7").Value = 10104401
$ws.Range("M7").Value = -3736.4
$ws.Range("N7").Value = -10104627

$ws.Range("H94").Value = 16146242
$ws.Range("I94").Value = 29412804
$ws.Range("J94").Value = 36846.355
$ws.Range("K94").Value = 29412804
$ws.Range("L94").Value = 36846.355
$ws.Range("M94").Value = -29412353
$ws.Range("N94").Value = -37748.355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1778.3
$ws.Range("I12").Value = 1947.7142
$ws.Range("J12").Value = 1383
$ws.Range("K12").Value = 1947.7142
$ws.Range("L12").Value = 1383
$ws.Range("M12").Value = -1777.7142
$ws.Range("N12").Value = -1723

$ws.Range("H105").Value = 4790.4287
$ws.Range("I105").Value = 3357.7
$ws.Range("J105").Value = 8372.25
$ws.Range("K105").Value = 3357.7
$ws.Range("L105").Value = 8372.25
$ws.Range("M105").Value = -1610.7
$ws.Range("N105").Value = -11866.25

$ws.Range("H132").Value = 1890.8096
$ws.Range("I132").Value = 1422.6666
$ws.Range("J132").Value = 4699.6665
$ws.Range("K132").Value = 4267.9998
$ws.Range("L132").Value = 14098.9995
$ws.Range("M132").Value = -1737.9998
$ws.Range("N132").Value = -19158.9995

$ws.Range("H134").Value = 1022.95654
$ws.Range("I134").Value = 891.44446
$ws.Range("J134").Value = 1496.4
$ws.Range("K134").Value = 2674.33338
$ws.Range("L134").Value = 4489.200000000001
$ws.Range("M134").Value = -139.33338
$ws.Range("N134").Value = -9559.200000000001

$ws.Range("H135").Value = 73333.336
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 73333.336
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 73333.336
$ws.Range("N135").Value = -83473.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 569.55554
$ws.Range("I23").Value = 70
$ws.Range("J23").Value = 632
$ws.Range("K23").Value = 210
$ws.Range("L23").Value = 1896
$ws.Range("M23").Value = 25
$ws.Range("N23").Value = -2366

$ws.Range("H32").Value = 9997
$ws.Range("I32").Value = 9995
$ws.Range("J32").Value = 9997.5
$ws.Range("K32").Value = 29985
$ws.Range("L32").Value = 29992.5
$ws.Range("M32").Value = -29702
$ws.Range("N32").Value = -30558.5

$ws.Range("H45").Value = 4466.6665
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4466.6665
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13399.9995
$ws.Range("N45").Value = -14463.9995

$ws.Range("H46").Value = 27481.21
$ws.Range("I46").Value = 1276.2106
$ws.Range("J46").Value = 53686.21
$ws.Range("K46").Value = 3828.6318
$ws.Range("L46").Value = 161058.63
$ws.Range("M46").Value = -3737.6318
$ws.Range("N46").Value = -161240.63

$ws.Range("H136").Value = 5809.4707
$ws.Range("I136").Value = 4780.136
$ws.Range("J136").Value = 7696.5835
$ws.Range("K136").Value = 14340.408
$ws.Range("L136").Value = 23089.7505
$ws.Range("M136").Value = -9240.408000000001
$ws.Range("N136").Value = -33289.75049999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3456660.5
$ws.Range("I3").Value = 5166081
$ws.Range("J3").Value = 2601950
$ws.Range("K3").Value = 5166081
$ws.Range("L3").Value = 2601950
$ws.Range("M3").Value = -5165965
$ws.Range("N3").Value = -2602182

$ws.Range("H107").Value = 1088.7059
$ws.Range("I107").Value = 1107.9286
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 1107.9286
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 812.0714
$ws.Range("N107").Value = -4839

$ws.Range("H132").Value = 4647.5
$ws.Range("I132").Value = 4402.45
$ws.Range("J132").Value = 6281.1665
$ws.Range("K132").Value = 13207.35
$ws.Range("L132").Value = 18843.4995
$ws.Range("M132").Value = -10677.35
$ws.Range("N132").Value = -23903.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8900
$ws.Range("I2").Value = 8430.77
$ws.Range("J2").Value = 15000
$ws.Range("K2").Value = 8430.77
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = -8318.77
$ws.Range("N2").Value = -15224

$ws.Range("H132").Value = 6971.048
$ws.Range("I132").Value = 7447.2104
$ws.Range("J132").Value = 2447.5
$ws.Range("K132").Value = 22341.6312
$ws.Range("L132").Value = 7342.5
$ws.Range("M132").Value = -19811.6312
$ws.Range("N132").Value = -12402.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 27749
$ws.Range("I2").Value = 27749
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 27749
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -27637

$ws.Range("H4").Value = 217929.7
$ws.Range("I4").Value = 184049.67
$ws.Range("J4").Value = 268749.75
$ws.Range("K4").Value = 184049.67
$ws.Range("L4").Value = 268749.75
$ws.Range("M4").Value = -183936.67
$ws.Range("N4").Value = -268975.75

$ws.Range("H5").Value = 231979.8
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 231979.8
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 231979.8
$ws.Range("N5").Value = -232203.8

$ws.Range("H96").Value = 2299.6667
$ws.Range("I96").Value = 2199
$ws.Range("J96").Value = 2350
$ws.Range("K96").Value = 2199
$ws.Range("L96").Value = 2350
$ws.Range("M96").Value = -826
$ws.Range("N96").Value = -5096

$ws.Range("H132").Value = 4621.3335
$ws.Range("I132").Value = 3275.4656
$ws.Range("J132").Value = 11717.728
$ws.Range("K132").Value = 9826.3968
$ws.Range("L132").Value = 35153.18399999999
$ws.Range("M132").Value = -7296.3968
$ws.Range("N132").Value = -40213.18399999999
